# Refresh the crypto price-tracker snapshot (Price / Volume(1h) columns)
# to the latest pulled values, per the "Updated symbol list" GitHub Action
# run. Only the Price (D) and Volume(1h) (E) columns move; row identity,
# coin name, link, date and hour columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> updated column values (only D and/or E are present per row,
# matching exactly the cells that changed in this snapshot).
$updates = @{
    2  = @{ "D" = "329.55";      "E" = "4.34%" };
    3  = @{ "D" = "40.60";       "E" = "7.49%" };
    4  = @{ "D" = "5.743";       "E" = "11.51%" };
    5  = @{ "D" = "0.08111";     "E" = "1.84%" };
    6  = @{ "D" = "4.598";       "E" = "2.87%" };
    7  = @{ "D" = "8.780";       "E" = "3.92%" };
    8  = @{ "D" = "1.970";       "E" = "2.78%" };
    9  = @{ "D" = "2.942";       "E" = "-1.19%" };
    10 = @{ "D" = "0.9450";      "E" = "0.51%" };
    11 = @{ "D" = "0.1285";      "E" = "1.47%" };
    12 = @{ "D" = "0.1987";      "E" = "2.94%" };
    13 = @{ "D" = "8.975";       "E" = "37.81%" };
    14 = @{ "D" = "0.09252";     "E" = "3.01%" };
    15 = @{ "D" = "0.03515";     "E" = "4.45%" };
    16 = @{ "D" = "0.09610";     "E" = "0.51%" };
    17 = @{ "D" = "0.001317";    "E" = "-3.43%" };
    18 = @{ "D" = "0.006107";    "E" = "0.73%" };
    19 = @{ "D" = "3.371";       "E" = "-1.06%" };
    20 = @{ "D" = "0.3566";      "E" = "1.58%" };
    21 = @{ "D" = "0.1425";      "E" = "9.62%" };
    22 = @{ "D" = "0.2411";      "E" = "4.72%" };
    23 = @{ "D" = "0.04430";     "E" = "2.18%" };
    24 = @{ "D" = "0.001258";    "E" = "4.99%" };
    25 = @{ "D" = "0.004358";    "E" = "-1.12%" };
    26 = @{ "D" = "0.0001091";   "E" = "-17.66%" };
    27 = @{                      "E" = "0.52%" };
    39 = @{ "D" = "0.02454";     "E" = "5.26%" };
    40 = @{ "D" = "0.05309";     "E" = "2.61%" };
    41 = @{ "D" = "0.007460";    "E" = "-0.48%" };
    42 = @{ "D" = "0.1430";      "E" = "2.31%" };
    43 = @{ "D" = "0.008724";    "E" = "1.49%" };
    44 = @{ "D" = "0.002103";    "E" = "5.53%" };
    45 = @{ "D" = "0.01088";     "E" = "37.41%" };
    46 = @{ "D" = "0.00006885";  "E" = "7.71%" };
    47 = @{ "D" = "0.00000000751"; "E" = "0.53%" };
    48 = @{ "D" = "0.003169";    "E" = "11.06%" };
    49 = @{                      "E" = "1.16%" };
    50 = @{ "D" = "0.00002102";  "E" = "0.53%" };
    51 = @{ "D" = "0.0002002";   "E" = "0.53%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $cell = $ws.Range($cellRef)
        # These columns are stored as plain text in the workbook (e.g. a
        # trailing-zero price like "40.60" or a percent string like
        # "4.34%"), so force text interpretation before writing - otherwise
        # Excel auto-converts the literal into a Number/Percentage and
        # drops formatting (trailing zeros, the literal "%" text, etc).
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        # Clear the one-off "Text" number format we just applied so the
        # cell's style matches its original (unstyled) state.
        $cell.Style = "Normal"
    }
}
